$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row for the 08/24/2025 bitcoin buy.
# The date column in this sheet is stored as plain text (e.g. "08/20/2025"
# in A42), not a real date. Assigning a date-looking string straight to
# .Value lets Excel's input-parsing auto-convert it into a date serial, so
# we briefly force the cell to Text format for the write, then clear the
# formatting again so the cell ends up unstyled/General, matching the rest
# of the column.
$dateCell = $ws.Range("A43")
$dateCell.NumberFormat = "@"
$dateCell.Value = "08/24/2025"
$dateCell.ClearFormats()

$ws.Range("B43").Value = 0.0004329499999999979
$ws.Range("C43").Value = 115486.7767640611
$ws.Range("D43").Value = 50
